$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "BA 450 Bachelor-Abschlussarbeit Bachelor Thesis" / "Scientific Work (SW)"
# row (row 83) is removed entirely; every row below it shifts up by one.
$ws.Rows(83).Delete()

# Reflect the updated view/selection state (user had scrolled back up and
# selected cell B5 after the edit).
$ws.Range("A4").Select()
$ws.Range("B5").Select()
